# YouTube-Comment-Analyzer-Setup.xlsx -- "Fix development mode - now functional"
#
# The Setup sheet's B1 cell held a sample YouTube URL that was wired up as a
# live hyperlink. This swaps in a different demo URL and removes the
# clickable hyperlink (the cell keeps looking like a hyperlink - blue,
# underlined - because it still carries the built-in "Hyperlink" cell
# style; only the clickable link / relationship goes away). The active
# selection is also left parked on the instructions' step-2 row (A7:B7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Swap the sample YouTube URL text in B1.
$ws.Range("B1").Value2 = "https://www.youtube.com/watch?v=MDE0outztmc&ab_channel=Curious%3F%3ANaturalWorld"

# 2) Remove the hyperlink (link + relationship) from the sheet, but keep the
#    cell's existing "Hyperlink" look/style (matches Excel's "Remove
#    Hyperlink" command, as opposed to "Clear Formats").
[void]$ws.Hyperlinks.Delete()

# 3) Leave the selection on A7:B7.
[void]$ws.Range("A7:B7").Select()
